# Add full names ("Team members list") to TEAM_MASTER sheet.
#
# The sheet has FirstName (D), MiddleName (E), LastName (F) split across
# three columns. This edit repurposes column D to hold the concatenated
# "First Middle Last" full name (MiddleName is blank for most rows, which
# is why most values below have a double space), while the original
# MiddleName/LastName columns (E/F) are left as-is.
#
# Row 14 (Samyukta Mazumder) additionally gets its LastName (F14) re-cased
# from the all-caps "MAZUMDER" to "Mazumder".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fullNames = @{
    2  = "Imran Khan Patan"
    3  = "Sudipta  Basak"
    4  = "Santosh  Kumar"
    5  = "Abhijit  Maiti"
    6  = "Sourav  Biswas"
    7  = "Alok Kumar Tripathi"
    8  = "Sumit  Das"
    9  = "Deepjyoti  Banerjee"
    10 = "Rajiv  Biswas"
    11 = "Raj  Kumar"
    12 = "Utsha   Rej"
    13 = "Pankaj Kumar Singh"
    14 = "Samyukta    Mazumder"
    15 = "Sajinur  Khatun"
    16 = "Vivek  Kumar"
    17 = "Anirban  Sardar"
    18 = "Akash  Bhattacharjee"
    19 = "Abhishek  Kaintura"
    20 = "Shaik  Alenoor"
    21 = "Preity  Mishra"
    22 = "Arpan  Halder"
    23 = "Sourav  Roy"
    24 = "Sreejita  Bose"
    25 = "Nabojita  Ghosh"
    26 = "Shreyashee  Majumder"
    27 = "Adarsh   Rana"
    28 = "Dwaipayan  Bhattacharyya"
    29 = "Sudip  Chowdhury"
    30 = "Rajarshi  Rakshit"
}

foreach ($row in $fullNames.Keys | Sort-Object) {
    $ws.Range("D$row").Value = $fullNames[$row]
}

# Re-case the LastName for the Samyukta Mazumder row to match the new
# full-name capitalization.
$ws.Range("F14").Value = "Mazumder"

# Column widths (best-fit) picked up by Excel once the longer full names
# are in place. ColumnWidth is specified in characters; Excel stores the
# rendered width with ~5/7 character of internal padding added on top, so
# we back that padding out here to land on the saved widths.
$padding = 5 / 7
$columnWidths = @{
    1 = 8
    2 = 8.5703125
    3 = 7.5703125
    4 = 21.85546875
    5 = 10.85546875
    6 = 12.140625
    7 = 27.140625
    8 = 6.85546875
    9 = 21.85546875
}
foreach ($col in $columnWidths.Keys | Sort-Object) {
    $ws.Columns.Item($col).ColumnWidth = $columnWidths[$col] - $padding
}

# Restore the selection to D4, matching the saved view state.
$ws.Range("D4").Select()
